$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name and title text to reflect the new date 2021-10-08
$ws.Name = "Through 2021-10-08"

# Update the "October (through 10-07)" label to "October (through 10-08)"
$ws.Range("A11").Value = "October (through 10-08)"

# Row 11: October monthly values
$ws.Range("C11").Value = 15
$ws.Range("D11").Value = 12
$ws.Range("E11").Value = 19
$ws.Range("F11").Value = 8
$ws.Range("H11").Value = 56

# Row 12: Total values
$ws.Range("C12").Value = 444
$ws.Range("D12").Value = 639
$ws.Range("E12").Value = 567
$ws.Range("F12").Value = 430
$ws.Range("H12").Value = 1305
